$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the percent style that previously lived on E2:E9 (it's being replaced
# with plain "N/A" text), but keep the Percent style on E10 since it will
# hold the averaged formula value.
$ws.Range("E2:E9").Style = "Normal"

# Fill the "N/A" placeholder across columns B:G for rows 2-9.
$ws.Range("B2:G9").Value = "N/A"

# Row 10 holds actual extracted data.
$ws.Range("B10").Value = 98200
$ws.Range("C10").Value = "N/A"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Formula = "=AVERAGE(0.33,0.51)"
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "N/A"

# Update the active selection to reflect where the user ended up after entry.
$ws.Range("J14").Select()
